$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2172284644194757
$ws.Range("C2").Value = 0.5205992509363296
$ws.Range("J2").Value = 0.0149812734082397
$ws.Range("P2").Value = 0.1385767790262172
$ws.Range("S2").Value = 0.1086142322097378

$ws.Range("B3").Value = 0.007042253521126761
$ws.Range("C3").Value = 0.02816901408450704
$ws.Range("J3").Value = 0.04225352112676056
$ws.Range("P3").Value = 0.7183098591549296
$ws.Range("S3").Value = 0.2042253521126761

$ws.Range("J4").Value = 0.02702702702702703
$ws.Range("P4").Value = 0.7567567567567568
$ws.Range("S4").Value = 0.2162162162162162

$ws.Range("B6").Value = 0.04145077720207254
$ws.Range("D6").Value = 0.02590673575129534
$ws.Range("F6").Value = 0.03626943005181347
$ws.Range("J6").Value = 0.2746113989637305
$ws.Range("O6").Value = 0.01036269430051814
$ws.Range("Q6").Value = 0.155440414507772
$ws.Range("R6").Value = 0.07253886010362694
$ws.Range("S6").Value = 0.383419689119171

$ws.Range("B7").Value = 0.09696969696969697
$ws.Range("D7").Value = 0.006060606060606061
$ws.Range("F7").Value = 0.03636363636363636
$ws.Range("J7").Value = 0.1393939393939394
$ws.Range("O7").Value = 0.006060606060606061
$ws.Range("Q7").Value = 0.1939393939393939
$ws.Range("R7").Value = 0.103030303030303
$ws.Range("S7").Value = 0.4181818181818182

$ws.Range("B8").Value = 0.09172259507829977
$ws.Range("D8").Value = 0.01342281879194631
$ws.Range("F8").Value = 0.04921700223713647
$ws.Range("J8").Value = 0.1140939597315436
$ws.Range("O8").Value = 0.01342281879194631
$ws.Range("Q8").Value = 0.1879194630872483
$ws.Range("R8").Value = 0.1006711409395973
$ws.Range("S8").Value = 0.4295302013422819

$ws.Range("B9").Value = 0.1004566210045662
$ws.Range("D9").Value = 0.0136986301369863
$ws.Range("F9").Value = 0.0684931506849315
$ws.Range("J9").Value = 0.0958904109589041
$ws.Range("O9").Value = 0.0091324200913242
$ws.Range("Q9").Value = 0.1780821917808219
$ws.Range("R9").Value = 0.1141552511415525
$ws.Range("S9").Value = 0.4200913242009132

$ws.Range("B10").Value = 0.09840201850294365
$ws.Range("D10").Value = 0.01934398654331371
$ws.Range("F10").Value = 0.0656013456686291
$ws.Range("J10").Value = 0.1295206055508831
$ws.Range("O10").Value = 0.01009251471825063
$ws.Range("Q10").Value = 0.208578637510513
$ws.Range("R10").Value = 0.0941968040370059
$ws.Range("S10").Value = 0.3742640874684609

$ws.Range("G11").Value = 0.1349206349206349
$ws.Range("J11").Value = 0.07936507936507936
$ws.Range("K11").Value = 0.1706349206349206
$ws.Range("L11").Value = 0.5952380952380952
$ws.Range("S11").Value = 0.01984126984126984

$ws.Range("G12").Value = 0.8079470198675497
$ws.Range("J12").Value = 0.1655629139072848
$ws.Range("K12").Value = 0.006622516556291391
$ws.Range("L12").Value = 0.01324503311258278
$ws.Range("S12").Value = 0.006622516556291391

$ws.Range("G13").Value = 0.5161290322580645
$ws.Range("J13").Value = 0.3870967741935484
$ws.Range("S13").Value = 0.0967741935483871

$ws.Range("G14").Value = 0.5
$ws.Range("J14").Value = 0.5

$ws.Range("F15").Value = 0.02702702702702703
$ws.Range("H15").Value = 0.1675675675675676
$ws.Range("I15").Value = 0.1027027027027027
$ws.Range("J15").Value = 0.3081081081081081
$ws.Range("K15").Value = 0.1027027027027027
$ws.Range("M15").Value = 0.02162162162162162
$ws.Range("O15").Value = 0.04864864864864865
$ws.Range("S15").Value = 0.2216216216216216

$ws.Range("F16").Value = 0.01219512195121951
$ws.Range("H16").Value = 0.1951219512195122
$ws.Range("I16").Value = 0.1097560975609756
$ws.Range("J16").Value = 0.3780487804878049
$ws.Range("K16").Value = 0.08536585365853659
$ws.Range("M16").Value = 0.01219512195121951
$ws.Range("O16").Value = 0.06097560975609756
$ws.Range("S16").Value = 0.1463414634146341

$ws.Range("F17").Value = 0.02102803738317757
$ws.Range("H17").Value = 0.1635514018691589
$ws.Range("I17").Value = 0.09345794392523364
$ws.Range("J17").Value = 0.4532710280373832
$ws.Range("K17").Value = 0.06775700934579439
$ws.Range("M17").Value = 0.01168224299065421
$ws.Range("N17").Value = 0.002336448598130841
$ws.Range("O17").Value = 0.07242990654205607
$ws.Range("S17").Value = 0.1144859813084112

$ws.Range("F18").Value = 0.01395348837209302
$ws.Range("H18").Value = 0.2
$ws.Range("I18").Value = 0.08372093023255814
$ws.Range("J18").Value = 0.4511627906976744
$ws.Range("K18").Value = 0.1116279069767442
$ws.Range("M18").Value = 0.0186046511627907
$ws.Range("O18").Value = 0.02325581395348837
$ws.Range("S18").Value = 0.09767441860465116

$ws.Range("F19").Value = 0.01571546732837055
$ws.Range("H19").Value = 0.2315963606286187
$ws.Range("I19").Value = 0.1025641025641026
$ws.Range("J19").Value = 0.347394540942928
$ws.Range("K19").Value = 0.0967741935483871
$ws.Range("M19").Value = 0.01654259718775848
$ws.Range("N19").Value = 0.0008271298593879239
$ws.Range("O19").Value = 0.06782464846980976
$ws.Range("S19").Value = 0.1207609594706369

